# Auto-generated Excel COM-interop script applying scheduled market-data refresh
# to the per-sheet profit tables (currentAveragePrice*, LevePrice*, LeveProfit*).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 42.5
$ws.Range("I8").Value = 42.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 127.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = 11.5

$ws.Range("H62").Value = 3999
$ws.Range("I62").Value = 3999
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3999
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3375
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 3999
$ws.Range("I65").Value = 3999
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 19995
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -16875
$ws.Range("N65").ClearContents()

$ws.Range("H92").Value = 898.8077
$ws.Range("I92").Value = 739.41174
$ws.Range("K92").Value = 739.41174
$ws.Range("M92").Value = 508.58826

$ws.Range("H137").Value = 2579.26
$ws.Range("I137").Value = 2502.122
$ws.Range("J137").Value = 2930.6667
$ws.Range("K137").Value = 7506.366
$ws.Range("L137").Value = 8792.000100000001
$ws.Range("M137").Value = -4956.366
$ws.Range("N137").Value = -13892.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 55575380
$ws.Range("I32").Value = 62521550
$ws.Range("K32").Value = 62521550
$ws.Range("M32").Value = -62521263

$ws.Range("H61").Value = 2832.625
$ws.Range("I61").Value = 2743.4614
$ws.Range("J61").Value = 3219
$ws.Range("K61").Value = 2743.4614
$ws.Range("L61").Value = 3219
$ws.Range("M61").Value = -2531.4614
$ws.Range("N61").Value = -3643

$ws.Range("H102").Value = 2013.4286
$ws.Range("I102").Value = 1682.3334
$ws.Range("K102").Value = 1682.3334
$ws.Range("M102").Value = -60.33339999999998

$ws.Range("H132").Value = 2072.6956
$ws.Range("I132").Value = 1985.091
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 5955.272999999999
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -3425.272999999999
$ws.Range("N132").Value = -17060

$ws.Range("H136").Value = 2832.625
$ws.Range("I136").Value = 2743.4614
$ws.Range("J136").Value = 3219
$ws.Range("K136").Value = 8230.3842
$ws.Range("L136").Value = 9657
$ws.Range("M136").Value = -5680.3842
$ws.Range("N136").Value = -14757

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2241.889
$ws.Range("I20").Value = 1407.0769
$ws.Range("K20").Value = 1407.0769
$ws.Range("M20").Value = -1160.0769

$ws.Range("H134").Value = 2634.3513
$ws.Range("I134").Value = 1879.7097
$ws.Range("J134").Value = 6533.3335
$ws.Range("K134").Value = 5639.1291
$ws.Range("L134").Value = 19600.0005
$ws.Range("M134").Value = -3104.1291
$ws.Range("N134").Value = -24670.0005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2082.1667
$ws.Range("I31").Value = 1951.1666
$ws.Range("J31").Value = 2475.1667
$ws.Range("K31").Value = 1951.1666
$ws.Range("L31").Value = 2475.1667
$ws.Range("M31").Value = -1656.1666
$ws.Range("N31").Value = -3065.1667

$ws.Range("H34").Value = 2082.1667
$ws.Range("I34").Value = 1951.1666
$ws.Range("J34").Value = 2475.1667
$ws.Range("K34").Value = 1951.1666
$ws.Range("L34").Value = 2475.1667
$ws.Range("M34").Value = -1749.1666
$ws.Range("N34").Value = -2879.1667

$ws.Range("H58").Value = 2493.3635
$ws.Range("I58").Value = 1727.0588
$ws.Range("K58").Value = 1727.0588
$ws.Range("M58").Value = -1524.0588

$ws.Range("H99").Value = 4962.706
$ws.Range("I99").Value = 2624.6
$ws.Range("J99").Value = 5936.9165
$ws.Range("K99").Value = 2624.6
$ws.Range("L99").Value = 5936.9165
$ws.Range("M99").Value = -1126.6
$ws.Range("N99").Value = -8932.916499999999

$ws.Range("H126").Value = 4962.706
$ws.Range("I126").Value = 2624.6
$ws.Range("J126").Value = 5936.9165
$ws.Range("K126").Value = 7873.799999999999
$ws.Range("L126").Value = 17810.7495
$ws.Range("M126").Value = -5403.799999999999
$ws.Range("N126").Value = -22750.7495

$ws.Range("H132").Value = 2175.0386
$ws.Range("I132").Value = 2208.3
$ws.Range("J132").Value = 2064.1667
$ws.Range("K132").Value = 6624.900000000001
$ws.Range("L132").Value = 6192.500100000001
$ws.Range("M132").Value = -4094.900000000001
$ws.Range("N132").Value = -11252.5001

$ws.Range("H134").Value = 7673.7144
$ws.Range("I134").Value = 9906
$ws.Range("K134").Value = 29718
$ws.Range("M134").Value = -27183

$ws.Range("H136").Value = 2493.3635
$ws.Range("I136").Value = 1727.0588
$ws.Range("K136").Value = 5181.1764
$ws.Range("M136").Value = -2631.1764

$ws.Range("H141").Value = 40199.668
$ws.Range("J141").Value = 40199.668
$ws.Range("L141").Value = 40199.668
$ws.Range("N141").Value = -50559.668

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 5149.619
$ws.Range("I97").Value = 343.70587
$ws.Range("J97").Value = 25574.75
$ws.Range("K97").Value = 343.70587
$ws.Range("L97").Value = 25574.75
$ws.Range("M97").Value = 152.29413
$ws.Range("N97").Value = -26566.75

$ws.Range("H113").Value = 5024.75
$ws.Range("I113").Value = 2266.3333
$ws.Range("K113").Value = 2266.3333
$ws.Range("M113").Value = -96.33329999999978

$ws.Range("H122").Value = 3835.342
$ws.Range("I122").Value = 2139.4443
$ws.Range("J122").Value = 5361.65
$ws.Range("K122").Value = 6418.3329
$ws.Range("L122").Value = 16084.95
$ws.Range("M122").Value = -3968.3329
$ws.Range("N122").Value = -20984.95

$ws.Range("H132").Value = 2645
$ws.Range("I132").Value = 2364.3333
$ws.Range("J132").Value = 6013
$ws.Range("K132").Value = 7092.999899999999
$ws.Range("L132").Value = 18039
$ws.Range("M132").Value = -4562.999899999999
$ws.Range("N132").Value = -23099

$ws.Range("H141").Value = 75995
$ws.Range("J141").Value = 75995
$ws.Range("L141").Value = 75995
$ws.Range("N141").Value = -86355

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4765.3887
$ws.Range("I7").Value = 3839.3333
$ws.Range("K7").Value = 3839.3333
$ws.Range("M7").Value = -3727.3333

$ws.Range("H40").Value = 10308.385
$ws.Range("I40").Value = 12110.4
$ws.Range("K40").Value = 12110.4
$ws.Range("M40").Value = -11974.4

$ws.Range("H74").Value = 76875
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 76875
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 76875
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -78871

$ws.Range("H77").Value = 76875
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 76875
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 230625
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -240609

$ws.Range("H93").Value = 4286
$ws.Range("I93").Value = 5581.2
$ws.Range("J93").Value = 2667
$ws.Range("K93").Value = 5581.2
$ws.Range("L93").Value = 2667
$ws.Range("M93").Value = -4333.2
$ws.Range("N93").Value = -5163

$ws.Range("H126").Value = 4765.3887
$ws.Range("I126").Value = 3839.3333
$ws.Range("K126").Value = 11517.9999
$ws.Range("M126").Value = -9047.999899999999

$ws.Range("H132").Value = 5102
$ws.Range("I132").Value = 3281.6428
$ws.Range("K132").Value = 9844.928400000001
$ws.Range("M132").Value = -7314.928400000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1061.0476
$ws.Range("I126").Value = 1107.2106
$ws.Range("J126").Value = 622.5
$ws.Range("K126").Value = 3321.6318
$ws.Range("L126").Value = 1867.5
$ws.Range("M126").Value = -851.6318000000001
$ws.Range("N126").Value = -6807.5

$ws.Range("H136").Value = 3154.3333
$ws.Range("I136").Value = 1931.5
$ws.Range("K136").Value = 5794.5
$ws.Range("M136").Value = -3244.5
